$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01821980476379394
$ws.Range("C2").Value = 0.03011903762817383
$ws.Range("D2").Value = 0.006804370880126953
$ws.Range("E2").Value = 0.01992926597595215
$ws.Range("F2").Value = 0.003295183181762695
$ws.Range("G2").Value = 0.06986837387084961
$ws.Range("H2").Value = 0.02082295417785645
$ws.Range("I2").Value = 0.02711653709411621
$ws.Range("J2").Value = 0.01193118095397949
$ws.Range("K2").Value = 0.02800393104553223
$ws.Range("L2").Value = 0.003326177597045898
$ws.Range("M2").Value = 0.02225141525268555
$ws.Range("B3").Value = 0.08473114967346192
$ws.Range("C3").Value = 0.02169523239135742
$ws.Range("D3").Value = 0.01310839653015137
$ws.Range("E3").Value = 0.008887815475463866
$ws.Range("F3").Value = 0.006555986404418945
$ws.Range("G3").Value = 0.006721735000610352
$ws.Range("H3").Value = 0.09351215362548829
$ws.Range("I3").Value = 0.02199573516845703
$ws.Range("J3").Value = 0.06733946800231934
$ws.Range("K3").Value = 0.02134394645690918
$ws.Range("L3").Value = 0.02733683586120605
$ws.Range("M3").Value = 0.01155438423156738
$ws.Range("B4").Value = 0.03235025405883789
$ws.Range("C4").Value = 0.008517026901245117
$ws.Range("D4").Value = 0.01539134979248047
$ws.Range("E4").Value = 0.00366363525390625
$ws.Range("F4").Value = 0.06216011047363281
$ws.Range("G4").Value = 0.009285259246826171
$ws.Range("H4").Value = 0.02718038558959961
$ws.Range("I4").Value = 0.01134395599365234
$ws.Range("J4").Value = 0.01680665016174316
$ws.Range("K4").Value = 0.0165278434753418
$ws.Range("L4").Value = 0.03703403472900391
$ws.Range("M4").Value = 0.006600427627563477
$ws.Range("B5").Value = 0.0191474437713623
$ws.Range("C5").Value = 0.01380510330200195
$ws.Range("D5").Value = 0.0155123233795166
$ws.Range("E5").Value = 0.01106185913085938
$ws.Range("H5").Value = 0.01542587280273438
$ws.Range("I5").Value = 0.01747221946716309
$ws.Range("J5").Value = 0.013922119140625
$ws.Range("K5").Value = 0.009642696380615235
$ws.Range("B6").Value = 0.3828513622283936
$ws.Range("C6").Value = 0.06287045478820801
$ws.Range("D6").Value = 0.2550122261047363
$ws.Range("E6").Value = 0.04659237861633301
$ws.Range("F6").Value = 0.1097476005554199
$ws.Range("G6").Value = 0.01874260902404785
$ws.Range("H6").Value = 0.3801012992858886
$ws.Range("I6").Value = 0.0534945011138916
$ws.Range("J6").Value = 0.2442757129669189
$ws.Range("K6").Value = 0.04940900802612305
$ws.Range("L6").Value = 0.1176186561584473
$ws.Range("M6").Value = 0.01742396354675293
